$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1944444444444444
$ws.Range("C2").Value = 0.5555555555555556
$ws.Range("P2").Value = 0.1388888888888889
$ws.Range("S2").Value = 0.1111111111111111
$ws.Range("J3").Value = 0.05263157894736842
$ws.Range("P3").Value = 0.7894736842105263
$ws.Range("S3").Value = 0.1578947368421053
$ws.Range("J4").Value = 0.1111111111111111
$ws.Range("P4").Value = 0.4444444444444444
$ws.Range("S4").Value = 0.4444444444444444
$ws.Range("B6").Value = 0.03846153846153846
$ws.Range("J6").Value = 0.5384615384615384
$ws.Range("O6").Value = 0.03846153846153846
$ws.Range("Q6").Value = 0.03846153846153846
$ws.Range("S6").Value = 0.3461538461538461
$ws.Range("F7").Value = 0.04347826086956522
$ws.Range("J7").Value = 0.2173913043478261
$ws.Range("Q7").Value = 0.2608695652173913
$ws.Range("R7").Value = 0.04347826086956522
$ws.Range("S7").Value = 0.4347826086956522
$ws.Range("B8").Value = 0.05263157894736842
$ws.Range("D8").Value = 0.02631578947368421
$ws.Range("F8").Value = 0.02631578947368421
$ws.Range("J8").Value = 0.2105263157894737
$ws.Range("O8").Value = 0.0131578947368421
$ws.Range("Q8").Value = 0.1842105263157895
$ws.Range("R8").Value = 0.1710526315789474
$ws.Range("S8").Value = 0.3157894736842105
$ws.Range("B9").Value = 0.06666666666666667
$ws.Range("J9").Value = 0.3333333333333333
$ws.Range("Q9").Value = 0.1333333333333333
$ws.Range("R9").Value = 0.2
$ws.Range("S9").Value = 0.2666666666666667
$ws.Range("B10").Value = 0.07954545454545454
$ws.Range("D10").Value = 0.02651515151515152
$ws.Range("F10").Value = 0.0303030303030303
$ws.Range("J10").Value = 0.1553030303030303
$ws.Range("O10").Value = 0.01136363636363636
$ws.Range("Q10").Value = 0.3181818181818182
$ws.Range("R10").Value = 0.1098484848484848
$ws.Range("S10").Value = 0.2689393939393939
$ws.Range("G11").Value = 0.1333333333333333
$ws.Range("J11").Value = 0.1333333333333333
$ws.Range("K11").Value = 0.2444444444444444
$ws.Range("L11").Value = 0.4888888888888889
$ws.Range("G12").Value = 0.6666666666666666
$ws.Range("J12").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.03703703703703703
$ws.Range("H15").Value = 0.1851851851851852
$ws.Range("J15").Value = 0.4074074074074074
$ws.Range("K15").Value = 0.1111111111111111
$ws.Range("O15").Value = 0.07407407407407407
$ws.Range("S15").Value = 0.1851851851851852
$ws.Range("H16").Value = 0.08333333333333333
$ws.Range("I16").Value = 0.08333333333333333
$ws.Range("J16").Value = 0.5416666666666666
$ws.Range("K16").Value = 0.04166666666666666
$ws.Range("M16").Value = 0.04166666666666666
$ws.Range("S16").Value = 0.2083333333333333
$ws.Range("F17").Value = 0.06666666666666667
$ws.Range("H17").Value = 0.2095238095238095
$ws.Range("I17").Value = 0.01904761904761905
$ws.Range("J17").Value = 0.4380952380952381
$ws.Range("K17").Value = 0.05714285714285714
$ws.Range("M17").Value = 0.02857142857142857
$ws.Range("O17").Value = 0.04761904761904762
$ws.Range("S17").Value = 0.1333333333333333
$ws.Range("F18").Value = 0.02173913043478261
$ws.Range("H18").Value = 0.08695652173913043
$ws.Range("I18").Value = 0.08695652173913043
$ws.Range("J18").Value = 0.5217391304347826
$ws.Range("K18").Value = 0.08695652173913043
$ws.Range("M18").Value = 0.02173913043478261
$ws.Range("S18").Value = 0.1739130434782609
$ws.Range("F19").Value = 0.01129943502824859
$ws.Range("H19").Value = 0.2485875706214689
$ws.Range("I19").Value = 0.03954802259887006
$ws.Range("J19").Value = 0.423728813559322
$ws.Range("K19").Value = 0.096045197740113
$ws.Range("M19").Value = 0.01129943502824859
$ws.Range("O19").Value = 0.06779661016949153
$ws.Range("S19").Value = 0.1016949152542373
